# Update "想去人数" (interest count) figures across the workbook.
# Sheet order in workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($SheetName, $Row, $Value)
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Cells.Item($Row, 6).Value = $Value
}

# 展览 sheet
Set-CellValue "展览" 6  846
Set-CellValue "展览" 8  1251
Set-CellValue "展览" 18 2986
Set-CellValue "展览" 23 318
Set-CellValue "展览" 26 5363
Set-CellValue "展览" 29 29

# 演出 sheet
Set-CellValue "演出" 26 3980

# 本地生活 sheet
Set-CellValue "本地生活" 5 2501
Set-CellValue "本地生活" 6 1063

# 全部类型 sheet
Set-CellValue "全部类型" 5  2501
Set-CellValue "全部类型" 7  1063
Set-CellValue "全部类型" 9  372
Set-CellValue "全部类型" 13 846
Set-CellValue "全部类型" 15 1251
Set-CellValue "全部类型" 25 2986
Set-CellValue "全部类型" 28 318
Set-CellValue "全部类型" 32 5363
Set-CellValue "全部类型" 37 29
